$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 53.0169945
$ws.Range("H2").Value = 106.033989
$ws.Range("I2").Value = 0.7442421144210264
$ws.Range("J2").Value = 0.7041603619966643
$ws.Range("M2").Value = 0.3250655
$ws.Range("N2").Value = 0.650131
$ws.Range("O2").Value = 0.03171104771407953
$ws.Range("P2").Value = 0.02364479350543662
$ws.Range("Q2").Value = 17.23399582563975
$ws.Range("R2").Value = 68.93598330255901
$ws.Range("S2").Value = 0.0236006972012326
$ws.Range("T2").Value = 0.01664972635412463
$ws.Range("G3").Value = 53.0169945
$ws.Range("H3").Value = 106.033989
$ws.Range("I3").Value = 0.7442421144210264
$ws.Range("J3").Value = 0.7041603619966643
$ws.Range("O3").Value = 0.6816872822276142
$ws.Range("P3").Value = 0.762432473166021
$ws.Range("Q3").Value = 370.476430871323
$ws.Range("R3").Value = 2222.858585227938
$ws.Range("S3").Value = 0.5073403842990025
$ws.Range("T3").Value = 0.5368747263025975
$ws.Range("G4").Value = 53.0169945
$ws.Range("H4").Value = 106.033989
$ws.Range("I4").Value = 0.7442421144210264
$ws.Range("J4").Value = 0.7041603619966643
$ws.Range("M4").Value = 2.9317775
$ws.Range("N4").Value = 5.863555
$ws.Range("O4").Value = 0.2860030861151516
$ws.Range("P4").Value = 0.2132532477035712
$ws.Range("Q4").Value = 155.4340315927238
$ws.Range("R4").Value = 621.736126370895
$ws.Range("S4").Value = 0.2128555415412793
$ws.Range("T4").Value = 0.150164484099911
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("G5").Value = 53.0169945
$ws.Range("H5").Value = 106.033989
$ws.Range("I5").Value = 0.7442421144210264
$ws.Range("J5").Value = 0.7041603619966643
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.006136
$ws.Range("N5").Value = 0.018408
$ws.Range("O5").Value = 0.0005985839431548165
$ws.Range("P5").Value = 0.0006694856249710864
$ws.Range("Q5").Value = 0.325312278252
$ws.Range("R5").Value = 1.951873669512
$ws.Range("S5").Value = 0.0004454913795120161
$ws.Range("T5").Value = 0.0004714252400312033
$ws.Range("I6").Value = 0.09075212572810222
$ws.Range("J6").Value = 0.1287968964663749
$ws.Range("M6").Value = 0.3250655
$ws.Range("N6").Value = 0.650131
$ws.Range("O6").Value = 0.03171104771407953
$ws.Range("P6").Value = 0.02364479350543662
$ws.Range("Q6").Value = 2.101495905244167
$ws.Range("R6").Value = 12.608975431465
$ws.Range("S6").Value = 0.002877844989117994
$ws.Range("T6").Value = 0.003045376021088535
$ws.Range("I7").Value = 0.09075212572810222
$ws.Range("J7").Value = 0.1287968964663749
$ws.Range("O7").Value = 0.6816872822276142
$ws.Range("P7").Value = 0.762432473166021
$ws.Range("S7").Value = 0.06186456994396874
$ws.Range("T7").Value = 0.09819893630896617
$ws.Range("I8").Value = 0.09075212572810222
$ws.Range("J8").Value = 0.1287968964663749
$ws.Range("M8").Value = 2.9317775
$ws.Range("N8").Value = 5.863555
$ws.Range("O8").Value = 0.2860030861151516
$ws.Range("P8").Value = 0.2132532477035712
$ws.Range("Q8").Value = 18.95346756680416
$ws.Range("R8").Value = 113.720805400825
$ws.Range("S8").Value = 0.02595538802974748
$ws.Range("T8").Value = 0.02746635646559506
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("I9").Value = 0.09075212572810222
$ws.Range("J9").Value = 0.1287968964663749
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.006136
$ws.Range("N9").Value = 0.018408
$ws.Range("O9").Value = 0.0005985839431548165
$ws.Range("P9").Value = 0.0006694856249710864
$ws.Range("Q9").Value = 0.03966824801333333
$ws.Range("R9").Value = 0.35701423212
$ws.Range("S9").Value = [double]"5.43227652680091E-05"
$ws.Range("T9").Value = [double]"8.622767072512732E-05"
$ws.Range("G10").Value = 0.6789063333333333
$ws.Range("H10").Value = 2.036719
$ws.Range("I10").Value = 0.00953035323444874
$ws.Range("J10").Value = 0.01352563269430035
$ws.Range("M10").Value = 0.3250655
$ws.Range("N10").Value = 0.650131
$ws.Range("O10").Value = 0.03171104771407953
$ws.Range("P10").Value = 0.02364479350543662
$ws.Range("Q10").Value = 0.2206890266981666
$ws.Range("R10").Value = 1.324134160189
$ws.Range("S10").Value = 0.0003022174861496362
$ws.Range("T10").Value = 0.0003198107920871142
$ws.Range("G11").Value = 0.6789063333333333
$ws.Range("H11").Value = 2.036719
$ws.Range("I11").Value = 0.00953035323444874
$ws.Range("J11").Value = 0.01352563269430035
$ws.Range("O11").Value = 0.6816872822276142
$ws.Range("P11").Value = 0.762432473166021
$ws.Range("Q11").Value = 4.744116441177555
$ws.Range("R11").Value = 42.69704797059799
$ws.Range("S11").Value = 0.006496720595060514
$ws.Range("T11").Value = 0.01031238158625061
$ws.Range("G12").Value = 0.6789063333333333
$ws.Range("H12").Value = 2.036719
$ws.Range("I12").Value = 0.00953035323444874
$ws.Range("J12").Value = 0.01352563269430035
$ws.Range("M12").Value = 2.9317775
$ws.Range("N12").Value = 5.863555
$ws.Range("O12").Value = 0.2860030861151516
$ws.Range("P12").Value = 0.2132532477035712
$ws.Range("Q12").Value = 1.990402312674167
$ws.Range("R12").Value = 11.942413876045
$ws.Range("S12").Value = 0.002725710436819856
$ws.Range("T12").Value = 0.002884385099305154
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("G13").Value = 0.6789063333333333
$ws.Range("H13").Value = 2.036719
$ws.Range("I13").Value = 0.00953035323444874
$ws.Range("J13").Value = 0.01352563269430035
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.006136
$ws.Range("N13").Value = 0.018408
$ws.Range("O13").Value = 0.0005985839431548165
$ws.Range("P13").Value = 0.0006694856249710864
$ws.Range("Q13").Value = 0.004165769261333333
$ws.Range("R13").Value = 0.037491923352
$ws.Range("S13").Value = [double]"5.704716418734586E-06"
$ws.Range("T13").Value = [double]"9.055216657473031E-06"
$ws.Range("G14").Value = 10.1095025
$ws.Range("H14").Value = 20.219005
$ws.Range("I14").Value = 0.1419152026119597
$ws.Range("J14").Value = 0.1342722462324073
$ws.Range("M14").Value = 0.3250655
$ws.Range("N14").Value = 0.650131
$ws.Range("O14").Value = 0.03171104771407953
$ws.Range("P14").Value = 0.02364479350543662
$ws.Range("Q14").Value = 3.28625048491375
$ws.Range("R14").Value = 13.145001939655
$ws.Range("S14").Value = 0.004500279761381117
$ws.Range("T14").Value = 0.003174839535676411
$ws.Range("G15").Value = 10.1095025
$ws.Range("H15").Value = 20.219005
$ws.Range("I15").Value = 0.1419152026119597
$ws.Range("J15").Value = 0.1342722462324073
$ws.Range("O15").Value = 0.6816872822276142
$ws.Range("P15").Value = 0.762432473166021
$ws.Range("Q15").Value = 70.64399706936833
$ws.Range("R15").Value = 423.86398241621
$ws.Range("S15").Value = 0.09674178877532801
$ws.Range("T15").Value = 0.1023735207725312
$ws.Range("G16").Value = 10.1095025
$ws.Range("H16").Value = 20.219005
$ws.Range("I16").Value = 0.1419152026119597
$ws.Range("J16").Value = 0.1342722462324073
$ws.Range("M16").Value = 2.9317775
$ws.Range("N16").Value = 5.863555
$ws.Range("O16").Value = 0.2860030861151516
$ws.Range("P16").Value = 0.2132532477035712
$ws.Range("Q16").Value = 29.63881196569375
$ws.Range("R16").Value = 118.555247862775
$ws.Range("S16").Value = 0.04058818591367749
$ws.Range("T16").Value = 0.02863399258551446
$ws.Range("D17").Value = "Resolving-Mac"
$ws.Range("G17").Value = 10.1095025
$ws.Range("H17").Value = 20.219005
$ws.Range("I17").Value = 0.1419152026119597
$ws.Range("J17").Value = 0.1342722462324073
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 0.3333333333333333
$ws.Range("M17").Value = 0.006136
$ws.Range("N17").Value = 0.018408
$ws.Range("O17").Value = 0.0005985839431548165
$ws.Range("P17").Value = 0.0006694856249710864
$ws.Range("Q17").Value = 0.06203190734
$ws.Range("R17").Value = 0.37219144404
$ws.Range("S17").Value = [double]"8.494816157308153E-05"
$ws.Range("T17").Value = [double]"8.989333868517481E-05"
$ws.Range("G18").Value = 0.5288903333333334
$ws.Range("H18").Value = 1.586671
$ws.Range("I18").Value = 0.007424458207958987
$ws.Range("J18").Value = 0.01053691213795238
$ws.Range("M18").Value = 0.3250655
$ws.Range("N18").Value = 0.650131
$ws.Range("O18").Value = 0.03171104771407953
$ws.Range("P18").Value = 0.02364479350543662
$ws.Range("Q18").Value = 0.1719240006501667
$ws.Range("R18").Value = 1.031544003901
$ws.Range("S18").Value = 0.0002354373484837768
$ws.Range("T18").Value = 0.0002491431116868129
$ws.Range("G19").Value = 0.5288903333333334
$ws.Range("H19").Value = 1.586671
$ws.Range("I19").Value = 0.007424458207958987
$ws.Range("J19").Value = 0.01053691213795238
$ws.Range("O19").Value = 0.6816872822276142
$ws.Range("P19").Value = 0.762432473166021
$ws.Range("Q19").Value = 3.695822535086889
$ws.Range("R19").Value = 33.262402815782
$ws.Range("S19").Value = 0.005061158737796065
$ws.Range("T19").Value = 0.008033683980872102
$ws.Range("G20").Value = 0.5288903333333334
$ws.Range("H20").Value = 1.586671
$ws.Range("I20").Value = 0.007424458207958987
$ws.Range("J20").Value = 0.01053691213795238
$ws.Range("M20").Value = 2.9317775
$ws.Range("N20").Value = 5.863555
$ws.Range("O20").Value = 0.2860030861151516
$ws.Range("P20").Value = 0.2132532477035712
$ws.Range("Q20").Value = 1.550588779234167
$ws.Range("R20").Value = 9.303532675405
$ws.Range("S20").Value = 0.002123417960209238
$ws.Range("T20").Value = 0.002247030734185525
$ws.Range("D21").Value = "Resolving-Mac"
$ws.Range("G21").Value = 0.5288903333333334
$ws.Range("H21").Value = 1.586671
$ws.Range("I21").Value = 0.007424458207958987
$ws.Range("J21").Value = 0.01053691213795238
$ws.Range("K21").Value = 1
$ws.Range("L21").Value = 0.3333333333333333
$ws.Range("M21").Value = 0.006136
$ws.Range("N21").Value = 0.018408
$ws.Range("O21").Value = 0.0005985839431548165
$ws.Range("P21").Value = 0.0006694856249710864
$ws.Range("Q21").Value = 0.003245271085333333
$ws.Range("R21").Value = 0.029207439768
$ws.Range("S21").Value = [double]"4.444161469908233E-06"
$ws.Range("T21").Value = [double]"7.054311207942478E-06"
$ws.Range("G22").Value = 0.4370873333333334
$ws.Range("H22").Value = 1.311262
$ws.Range("I22").Value = 0.006135745796503949
$ws.Range("J22").Value = 0.008707950472300633
$ws.Range("M22").Value = 0.3250655
$ws.Range("N22").Value = 0.650131
$ws.Range("O22").Value = 0.03171104771407953
$ws.Range("P22").Value = 0.02364479350543662
$ws.Range("Q22").Value = 0.1420820125536667
$ws.Range("R22").Value = 0.8524920753220001
$ws.Range("S22").Value = 0.0001945709277143997
$ws.Range("T22").Value = 0.0002058976907731178
$ws.Range("G23").Value = 0.4370873333333334
$ws.Range("H23").Value = 1.311262
$ws.Range("I23").Value = 0.006135745796503949
$ws.Range("J23").Value = 0.008707950472300633
$ws.Range("O23").Value = 0.6816872822276142
$ws.Range("P23").Value = 0.762432473166021
$ws.Range("Q23").Value = 3.054314126244889
$ws.Range("R23").Value = 27.488827136204
$ws.Range("S23").Value = 0.004182659876458285
$ws.Range("T23").Value = 0.006639224214803393
$ws.Range("G24").Value = 0.4370873333333334
$ws.Range("H24").Value = 1.311262
$ws.Range("I24").Value = 0.006135745796503949
$ws.Range("J24").Value = 0.008707950472300633
$ws.Range("M24").Value = 2.9317775
$ws.Range("N24").Value = 5.863555
$ws.Range("O24").Value = 0.2860030861151516
$ws.Range("P24").Value = 0.2132532477035712
$ws.Range("Q24").Value = 1.281442809401667
$ws.Range("R24").Value = 7.688656856410001
$ws.Range("S24").Value = 0.001754842233418198
$ws.Range("T24").Value = 0.001856998719059957
$ws.Range("D25").Value = "Resolving-Mac"
$ws.Range("G25").Value = 0.4370873333333334
$ws.Range("H25").Value = 1.311262
$ws.Range("I25").Value = 0.006135745796503949
$ws.Range("J25").Value = 0.008707950472300633
$ws.Range("K25").Value = 1
$ws.Range("L25").Value = 0.3333333333333333
$ws.Range("M25").Value = 0.006136
$ws.Range("N25").Value = 0.018408
$ws.Range("O25").Value = 0.0005985839431548165
$ws.Range("P25").Value = 0.0006694856249710864
$ws.Range("Q25").Value = 0.002681967877333334
$ws.Range("R25").Value = 0.024137710896
$ws.Range("S25").Value = [double]"3.672758913066924E-06"
$ws.Range("T25").Value = [double]"5.829847664165457E-06"
